# Generate Report for Handoff
# Updates the "f889c0fd-9ec9-4624-b0c8-bc6f7f90dd34" row across the Overview,
# zh-cn and de-de sheets to reflect that the item is ready for handoff again
# (the previously handed-back file is stale vs. the latest source).

$wb = $excel.ActiveWorkbook

$readyForHandoff = "Ready for handoff"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/081cc6427214470329dfc2c373002892e2dad66f/e2e/f889c0fd-9ec9-4624-b0c8-bc6f7f90dd34.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d6c2866334750cb1e3cee4950982ad2eaf224e6d/e2e/f889c0fd-9ec9-4624-b0c8-bc6f7f90dd34.md."

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $readyForHandoff
$wsOverview.Range("F3").Value = $readyForHandoff
$wsOverview.Range("G3").Value = "2016-08-23 18:51:46"

# ---- zh-cn sheet ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $readyForHandoff
$wsZhCn.Range("H3").Value = "2016-08-23 18:51:41"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# ---- de-de sheet ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $readyForHandoff
$wsDeDe.Range("H3").Value = "2016-08-23 18:51:46"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
